$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.964396016490747
$ws.Range("D2").Value = 4.206718196610706
$ws.Range("E2").Value = 16.51181434500207
$ws.Range("F2").Value = 24.86926765575648
$ws.Range("G2").Value = 30.63217645484141
$ws.Range("H2").Value = 14.09805022422436
$ws.Range("I2").Value = 18.78269164326545
$ws.Range("K2").Value = 11.98120274892911
$ws.Range("N2").Value = 17.45144517107364
$ws.Range("B3").Value = 5.863659111749999
$ws.Range("D3").Value = 4.212820874785823
$ws.Range("E3").Value = 15.5719435487692
$ws.Range("F3").Value = 24.67859599536554
$ws.Range("G3").Value = 30.19873842431279
$ws.Range("H3").Value = 14.10734566929737
$ws.Range("I3").Value = 18.87901919925398
$ws.Range("K3").Value = 11.47004450851071
$ws.Range("N3").Value = 17.51823653641463
$ws.Range("B4").Value = 5.801872831340034
$ws.Range("D4").Value = 4.216732836313295
$ws.Range("E4").Value = 14.97027272857293
$ws.Range("F4").Value = 24.57000606805609
$ws.Range("G4").Value = 29.94318926086584
$ws.Range("H4").Value = 14.11673086142441
$ws.Range("I4").Value = 18.9413650936841
$ws.Range("K4").Value = 11.14686325558496
$ws.Range("N4").Value = 17.56103399800652
$ws.Range("B5").Value = 5.776744075375455
$ws.Range("D5").Value = 4.218368799224925
$ws.Range("E5").Value = 14.71917440459052
$ws.Range("F5").Value = 24.52792788309236
$ws.Range("G5").Value = 29.84185166896122
$ws.Range("H5").Value = 14.1214772053647
$ws.Range("I5").Value = 18.96757746869731
$ws.Range("K5").Value = 11.01302369953381
$ws.Range("N5").Value = 17.5789254519833
$ws.Range("B6").Value = 5.772575433372279
$ws.Range("D6").Value = 4.218642984868145
$ws.Range("E6").Value = 14.67713155076216
$ws.Range("F6").Value = 24.52107315356945
$ws.Range("G6").Value = 29.82519771984219
$ws.Range("H6").Value = 14.12232090998887
$ws.Range("I6").Value = 18.97197872345448
$ws.Range("K6").Value = 10.99067694978259
$ws.Range("N6").Value = 17.58192361228029
$ws.Range("B7").Value = 5.801533691298077
$ws.Range("D7").Value = 4.216754729786631
$ws.Range("E7").Value = 14.96690987275108
$ws.Range("F7").Value = 24.56942973918401
$ws.Range("G7").Value = 29.94181106819096
$ws.Range("H7").Value = 14.11679114434066
$ws.Range("I7").Value = 18.94171533830383
$ws.Range("K7").Value = 11.14506662604725
$ws.Range("N7").Value = 17.56127345953503
$ws.Range("B8").Value = 5.929667712586024
$ws.Range("D8").Value = 4.208788380434776
$ws.Range("E8").Value = 16.19299596485927
$ws.Range("F8").Value = 24.80178876902169
$ws.Range("G8").Value = 30.480616359459
$ws.Range("H8").Value = 14.10049030809214
$ws.Range("I8").Value = 18.8152418374671
$ws.Range("K8").Value = 11.80700668655098
$ws.Range("N8").Value = 17.47410491711418
$ws.Range("B9").Value = 6.180048273255006
$ws.Range("D9").Value = 4.194459721112154
$ws.Range("E9").Value = 18.48579666097279
$ws.Range("F9").Value = 25.32277469928891
$ws.Range("G9").Value = 31.61471815529229
$ws.Range("H9").Value = 14.09782052512199
$ws.Range("I9").Value = 18.59256814062487
$ws.Range("K9").Value = 13.04027989465743
$ws.Range("N9").Value = 17.31726956748803
$ws.Range("B10").Value = 6.361616370626757
$ws.Range("D10").Value = 4.184700796682394
$ws.Range("E10").Value = 20.13690233936268
$ws.Range("F10").Value = 25.74242430399479
$ws.Range("G10").Value = 32.4861849693113
$ws.Range("H10").Value = 14.11384711547207
$ws.Range("I10").Value = 18.44434298452655
$ws.Range("K10").Value = 14.06751450193758
$ws.Range("N10").Value = 17.21052629457299
$ws.Range("B11").Value = 6.44332819885941
$ws.Range("D11").Value = 4.18042396008106
$ws.Range("E11").Value = 20.84606385895496
$ws.Range("F11").Value = 25.94064017725295
$ws.Range("G11").Value = 32.88897824281459
$ws.Range("H11").Value = 14.12506089348186
$ws.Range("I11").Value = 18.38023435703142
$ws.Range("K11").Value = 14.50923658240671
$ws.Range("N11").Value = 17.16378381510475
$ws.Range("B12").Value = 6.47411401856504
$ws.Range("D12").Value = 4.178827500822298
$ws.Range("E12").Value = 21.10861296275977
$ws.Range("F12").Value = 26.01668819234943
$ws.Range("G12").Value = 33.04225609255004
$ws.Range("H12").Value = 14.12987184832776
$ws.Range("I12").Value = 18.35643446332547
$ws.Range("K12").Value = 14.67282427146706
$ws.Range("N12").Value = 17.14634289736042
$ws.Range("B13").Value = 6.467491164008046
$ws.Range("D13").Value = 4.1791703045921
$ws.Range("E13").Value = 21.05233438544342
$ws.Range("F13").Value = 26.00026699218747
$ws.Range("G13").Value = 33.00921427509622
$ws.Range("H13").Value = 14.1288106147066
$ws.Range("I13").Value = 18.3615390080045
$ws.Range("K13").Value = 14.63775657752583
$ws.Range("N13").Value = 17.15008759810757
$ws.Range("B14").Value = 6.445864270637899
$ws.Range("D14").Value = 4.180292157353406
$ws.Range("E14").Value = 20.86778387903248
$ws.Range("F14").Value = 25.94687726677566
$ws.Range("G14").Value = 32.90157443760876
$ws.Range("H14").Value = 14.12544537914627
$ws.Range("I14").Value = 18.37826677659873
$ws.Range("K14").Value = 14.52276886758854
$ws.Range("N14").Value = 17.16234374998077
$ws.Range("B15").Value = 6.432595919186944
$ws.Range("D15").Value = 4.18098232266713
$ws.Range("E15").Value = 20.75396168375334
$ws.Range("F15").Value = 25.91430124604475
$ws.Range("G15").Value = 32.83573456510347
$ws.Range("H15").Value = 14.12345759711354
$ws.Range("I15").Value = 18.38857507144227
$ws.Range("K15").Value = 14.45185586143742
$ws.Range("N15").Value = 17.16988473589308
$ws.Range("B16").Value = 6.356255959641417
$ws.Range("D16").Value = 4.184983543987319
$ws.Range("E16").Value = 20.08971614493692
$ws.Range("F16").Value = 25.7296122149674
$ws.Range("G16").Value = 32.45997523355138
$ws.Range("H16").Value = 14.11319325941311
$ws.Range("I16").Value = 18.44859937288861
$ws.Range("K16").Value = 14.03813176708431
$ws.Range("N16").Value = 17.21361741983759
$ws.Range("B17").Value = 6.309175910270202
$ws.Range("D17").Value = 4.187479587286276
$ws.Range("E17").Value = 19.67150763302406
$ws.Range("F17").Value = 25.61814029591926
$ws.Range("G17").Value = 32.23096926182462
$ws.Range("H17").Value = 14.10790169866105
$ws.Range("I17").Value = 18.48627207313397
$ws.Range("K17").Value = 13.77777037675449
$ws.Range("N17").Value = 17.24090985337526
$ws.Range("B18").Value = 6.282015143371317
$ws.Range("D18").Value = 4.188930561607215
$ws.Range("E18").Value = 19.42701597012315
$ws.Range("F18").Value = 25.55471712224052
$ws.Range("G18").Value = 32.09985942380217
$ws.Range("H18").Value = 14.10522741502943
$ws.Range("I18").Value = 18.50825292376232
$ws.Range("K18").Value = 13.62561045009793
$ws.Range("N18").Value = 17.25677872673184
$ws.Range("B19").Value = 6.272805887178967
$ws.Range("D19").Value = 4.189424476103628
$ws.Range("E19").Value = 19.34355594212839
$ws.Range("F19").Value = 25.53336397855948
$ws.Range("G19").Value = 32.05557744490047
$ws.Range("H19").Value = 14.1043853521439
$ws.Range("I19").Value = 18.51574895872655
$ws.Range("K19").Value = 13.5736785633743
$ws.Range("N19").Value = 17.26218107161356
$ws.Range("B20").Value = 6.314196322066776
$ws.Range("D20").Value = 4.187212296230896
$ws.Range("E20").Value = 19.71643511249298
$ws.Range("F20").Value = 25.62993545470038
$ws.Range("G20").Value = 32.25528558395074
$ws.Range("H20").Value = 14.10842677064631
$ws.Range("I20").Value = 18.48222941180163
$ws.Range("K20").Value = 13.80573545485052
$ws.Range("N20").Value = 17.23798684261624
$ws.Range("B21").Value = 6.452221093551957
$ws.Range("D21").Value = 4.179962017695664
$ws.Range("E21").Value = 20.92215325870183
$ws.Range("F21").Value = 25.96253283038164
$ws.Range("G21").Value = 32.93317185397156
$ws.Range("H21").Value = 14.1264185085781
$ws.Range("I21").Value = 18.37334049478319
$ws.Range("K21").Value = 14.5566435214455
$ws.Range("N21").Value = 17.15873679280742
$ws.Range("B22").Value = 6.54150238262696
$ws.Range("D22").Value = 4.175357978318597
$ws.Range("E22").Value = 21.67523267020787
$ws.Range("F22").Value = 26.18563201188678
$ws.Range("G22").Value = 33.3804992374948
$ws.Range("H22").Value = 14.14146751236327
$ws.Range("I22").Value = 18.30495328642645
$ws.Range("K22").Value = 15.02593894926466
$ws.Range("N22").Value = 17.10845379246705
$ws.Range("B23").Value = 6.493945365385303
$ws.Range("D23").Value = 4.177803032094714
$ws.Range("E23").Value = 21.27648427110521
$ws.Range("F23").Value = 26.06605735741926
$ws.Range("G23").Value = 33.1414137641979
$ws.Range("H23").Value = 14.13313452201195
$ws.Range("I23").Value = 18.3411988658052
$ws.Range("K23").Value = 14.77743212570614
$ws.Range("N23").Value = 17.1351530128361
$ws.Range("B24").Value = 6.311926883872894
$ws.Range("D24").Value = 4.187333088731136
$ws.Range("E24").Value = 19.69613604198521
$ws.Range("F24").Value = 25.6246007917314
$ws.Range("G24").Value = 32.24429045273982
$ws.Range("H24").Value = 14.10818823960145
$ws.Range("I24").Value = 18.48405609494242
$ws.Range("K24").Value = 13.79310015080069
$ws.Range("N24").Value = 17.23930778089778
$ws.Range("B25").Value = 6.112603153591444
$ws.Range("D25").Value = 4.198199662721435
$ws.Range("E25").Value = 17.84034138116694
$ws.Range("F25").Value = 25.17512922582336
$ws.Range("G25").Value = 31.30053354251994
$ws.Range("H25").Value = 14.09539113763984
$ws.Range("I25").Value = 18.65010232552279
$ws.Range("K25").Value = 12.7039989274611
$ws.Range("N25").Value = 17.35819970451882
